$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 2).Value = 0.301
$ws.Cells.Item(4, 3).Value = 0.053
$ws.Cells.Item(4, 4).Value = 0.231
$ws.Cells.Item(4, 5).Value = 0.162
$ws.Cells.Item(4, 7).Value = 0.111
$ws.Cells.Item(4, 8).Value = 0.204
$ws.Cells.Item(4, 10).Value = 0.098
$ws.Cells.Item(4, 11).Value = 0.382
$ws.Cells.Item(4, 12).Value = 0.099
$ws.Cells.Item(4, 13).Value = 0.314
$ws.Cells.Item(4, 14).Value = 0.286
$ws.Cells.Item(4, 15).Value = 0.018
$ws.Cells.Item(4, 16).Value = 0.134
$ws.Cells.Item(4, 17).Value = 0.574
$ws.Cells.Item(4, 18).Value = 0.208
$ws.Cells.Item(4, 19).Value = 0.456
$ws.Cells.Item(4, 20).Value = 0.316
$ws.Cells.Item(4, 22).Value = 0.297
$ws.Cells.Item(4, 23).Value = 0.263
$ws.Cells.Item(4, 25).Value = 0.206
$ws.Cells.Item(4, 26).Value = 0.467
$ws.Cells.Item(4, 27).Value = 0.126
$ws.Cells.Item(4, 28).Value = 0.355
$ws.Cells.Item(4, 31).Value = 0.07199999999999999
$ws.Cells.Item(4, 32).Value = 0.728
$ws.Cells.Item(4, 33).Value = 0.092
$ws.Cells.Item(4, 34).Value = 0.303
$ws.Cells.Item(4, 35).Value = 0.6909999999999999
$ws.Cells.Item(4, 36).Value = 0.162
$ws.Cells.Item(4, 37).Value = 0.402
$ws.Cells.Item(4, 38).Value = 0.746
$ws.Cells.Item(4, 39).Value = 0.105
$ws.Cells.Item(4, 40).Value = 0.325
$ws.Cells.Item(4, 41).Value = 0.722
$ws.Cells.Item(5, 2).Value = 0.824
$ws.Cells.Item(5, 3).Value = 0.145
$ws.Cells.Item(5, 4).Value = 0.381
$ws.Cells.Item(5, 5).Value = 0.706
$ws.Cells.Item(5, 6).Value = 0.208
$ws.Cells.Item(5, 7).Value = 0.456
$ws.Cells.Item(5, 8).Value = 0.882
$ws.Cells.Item(5, 9).Value = 0.104
$ws.Cells.Item(5, 10).Value = 0.322
$ws.Cells.Item(5, 11).Value = 0.735
$ws.Cells.Item(5, 12).Value = 0.195
$ws.Cells.Item(5, 13).Value = 0.441
$ws.Cells.Item(5, 14).Value = 0.882
$ws.Cells.Item(5, 15).Value = 0.104
$ws.Cells.Item(5, 16).Value = 0.322
$ws.Cells.Item(5, 17).Value = 0.647
$ws.Cells.Item(5, 18).Value = 0.228
$ws.Cells.Item(5, 19).Value = 0.478
$ws.Cells.Item(5, 20).Value = 0.647
$ws.Cells.Item(5, 21).Value = 0.228
$ws.Cells.Item(5, 22).Value = 0.478
$ws.Cells.Item(5, 23).Value = 0.794
$ws.Cells.Item(5, 24).Value = 0.163
$ws.Cells.Item(5, 25).Value = 0.404
$ws.Cells.Item(5, 26).Value = 0.853
$ws.Cells.Item(5, 27).Value = 0.125
$ws.Cells.Item(5, 28).Value = 0.354
$ws.Cells.Item(5, 29).Value = 0.824
$ws.Cells.Item(5, 30).Value = 0.145
$ws.Cells.Item(5, 31).Value = 0.381
$ws.Cells.Item(5, 32).Value = 0.971
$ws.Cells.Item(5, 34).Value = 0.169
$ws.Cells.Item(5, 35).Value = 0.794
$ws.Cells.Item(5, 36).Value = 0.163
$ws.Cells.Item(5, 37).Value = 0.404
$ws.Cells.Item(5, 38).Value = 0.9409999999999999
$ws.Cells.Item(5, 39).Value = 0.055
$ws.Cells.Item(5, 40).Value = 0.235
$ws.Cells.Item(5, 41).Value = 0.902
$ws.Cells.Item(6, 2).Value = 0.441
$ws.Cells.Item(6, 5).Value = 0.264
$ws.Cells.Item(6, 8).Value = 0.331
$ws.Cells.Item(6, 11).Value = 0.503
$ws.Cells.Item(6, 14).Value = 0.432
$ws.Cells.Item(6, 17).Value = 0.608
$ws.Cells.Item(6, 20).Value = 0.425
$ws.Cells.Item(6, 23).Value = 0.395
$ws.Cells.Item(6, 26).Value = 0.604
$ws.Cells.Item(6, 32).Value = 0.832
$ws.Cells.Item(6, 35).Value = 0.739
$ws.Cells.Item(6, 38).Value = 0.832
$ws.Cells.Item(6, 41).Value = 0.801
$ws.Cells.Item(7, 2).Value = 0.611
$ws.Cells.Item(7, 5).Value = 0.422
$ws.Cells.Item(7, 8).Value = 0.53
$ws.Cells.Item(7, 11).Value = 0.62
$ws.Cells.Item(7, 14).Value = 0.623
$ws.Cells.Item(7, 17).Value = 0.631
$ws.Cells.Item(7, 20).Value = 0.535
$ws.Cells.Item(7, 23).Value = 0.5659999999999999
$ws.Cells.Item(7, 26).Value = 0.732
$ws.Cells.Item(7, 29).Value = 0.408
$ws.Cells.Item(7, 32).Value = 0.91
$ws.Cells.Item(7, 35).Value = 0.771
$ws.Cells.Item(7, 38).Value = 0.894
$ws.Cells.Item(7, 41).Value = 0.858
$ws.Cells.Item(8, 2).Value = 0.762
$ws.Cells.Item(8, 3).Value = 0.146
$ws.Cells.Item(8, 4).Value = 0.383
$ws.Cells.Item(8, 5).Value = 0.594
$ws.Cells.Item(8, 6).Value = 0.184
$ws.Cells.Item(8, 7).Value = 0.429
$ws.Cells.Item(8, 8).Value = 0.768
$ws.Cells.Item(8, 9).Value = 0.122
$ws.Cells.Item(8, 10).Value = 0.349
$ws.Cells.Item(8, 11).Value = 0.655
$ws.Cells.Item(8, 12).Value = 0.183
$ws.Cells.Item(8, 13).Value = 0.428
$ws.Cells.Item(8, 14).Value = 0.793
$ws.Cells.Item(8, 15).Value = 0.115
$ws.Cells.Item(8, 16).Value = 0.339
$ws.Cells.Item(8, 17).Value = 0.614
$ws.Cells.Item(8, 18).Value = 0.216
$ws.Cells.Item(8, 19).Value = 0.465
$ws.Cells.Item(8, 20).Value = 0.5600000000000001
$ws.Cells.Item(8, 21).Value = 0.197
$ws.Cells.Item(8, 22).Value = 0.444
$ws.Cells.Item(8, 23).Value = 0.71
$ws.Cells.Item(8, 24).Value = 0.157
$ws.Cells.Item(8, 25).Value = 0.396
$ws.Cells.Item(8, 26).Value = 0.791
$ws.Cells.Item(8, 27).Value = 0.13
$ws.Cells.Item(8, 28).Value = 0.361
$ws.Cells.Item(8, 29).Value = 0.702
$ws.Cells.Item(8, 30).Value = 0.154
$ws.Cells.Item(8, 31).Value = 0.392
$ws.Cells.Item(8, 32).Value = 0.891
$ws.Cells.Item(8, 33).Value = 0.049
$ws.Cells.Item(8, 34).Value = 0.221
$ws.Cells.Item(8, 35).Value = 0.783
$ws.Cells.Item(8, 36).Value = 0.163
$ws.Cells.Item(8, 37).Value = 0.404
$ws.Cells.Item(8, 38).Value = 0.909
$ws.Cells.Item(8, 39).Value = 0.062
$ws.Cells.Item(8, 40).Value = 0.25
$ws.Cells.Item(8, 41).Value = 0.861
$ws.Cells.Item(9, 2).Value = 0.676
$ws.Cells.Item(9, 3).Value = 0.219
$ws.Cells.Item(9, 4).Value = 0.468
$ws.Cells.Item(9, 5).Value = 0.471
$ws.Cells.Item(9, 6).Value = 0.249
$ws.Cells.Item(9, 7).Value = 0.499
$ws.Cells.Item(9, 8).Value = 0.647
$ws.Cells.Item(9, 9).Value = 0.228
$ws.Cells.Item(9, 10).Value = 0.478
$ws.Cells.Item(9, 11).Value = 0.5590000000000001
$ws.Cells.Item(9, 12).Value = 0.247
$ws.Cells.Item(9, 13).Value = 0.497
$ws.Cells.Item(9, 14).Value = 0.676
$ws.Cells.Item(9, 15).Value = 0.219
$ws.Cells.Item(9, 16).Value = 0.468
$ws.Cells.Item(9, 17).Value = 0.5590000000000001
$ws.Cells.Item(9, 18).Value = 0.247
$ws.Cells.Item(9, 19).Value = 0.497
$ws.Cells.Item(9, 20).Value = 0.441
$ws.Cells.Item(9, 21).Value = 0.247
$ws.Cells.Item(9, 22).Value = 0.497
$ws.Cells.Item(9, 23).Value = 0.588
$ws.Cells.Item(9, 24).Value = 0.242
$ws.Cells.Item(9, 25).Value = 0.492
$ws.Cells.Item(9, 26).Value = 0.706
$ws.Cells.Item(9, 27).Value = 0.208
$ws.Cells.Item(9, 28).Value = 0.456
$ws.Cells.Item(9, 29).Value = 0.588
$ws.Cells.Item(9, 30).Value = 0.242
$ws.Cells.Item(9, 31).Value = 0.492
$ws.Cells.Item(9, 32).Value = 0.765
$ws.Cells.Item(9, 33).Value = 0.18
$ws.Cells.Item(9, 34).Value = 0.424
$ws.Cells.Item(9, 35).Value = 0.765
$ws.Cells.Item(9, 36).Value = 0.18
$ws.Cells.Item(9, 37).Value = 0.424
$ws.Cells.Item(9, 38).Value = 0.853
$ws.Cells.Item(9, 39).Value = 0.125
$ws.Cells.Item(9, 40).Value = 0.354
$ws.Cells.Item(9, 41).Value = 0.794
$ws.Cells.Item(10, 2).Value = 0.824
$ws.Cells.Item(10, 3).Value = 0.145
$ws.Cells.Item(10, 4).Value = 0.381
$ws.Cells.Item(10, 5).Value = 0.618
$ws.Cells.Item(10, 6).Value = 0.236
$ws.Cells.Item(10, 7).Value = 0.486
$ws.Cells.Item(10, 8).Value = 0.794
$ws.Cells.Item(10, 9).Value = 0.163
$ws.Cells.Item(10, 10).Value = 0.404
$ws.Cells.Item(10, 11).Value = 0.735
$ws.Cells.Item(10, 12).Value = 0.195
$ws.Cells.Item(10, 13).Value = 0.441
$ws.Cells.Item(10, 14).Value = 0.853
$ws.Cells.Item(10, 15).Value = 0.125
$ws.Cells.Item(10, 16).Value = 0.354
$ws.Cells.Item(10, 17).Value = 0.647
$ws.Cells.Item(10, 18).Value = 0.228
$ws.Cells.Item(10, 19).Value = 0.478
$ws.Cells.Item(10, 20).Value = 0.647
$ws.Cells.Item(10, 21).Value = 0.228
$ws.Cells.Item(10, 22).Value = 0.478
$ws.Cells.Item(10, 23).Value = 0.794
$ws.Cells.Item(10, 24).Value = 0.163
$ws.Cells.Item(10, 25).Value = 0.404
$ws.Cells.Item(10, 26).Value = 0.853
$ws.Cells.Item(10, 27).Value = 0.125
$ws.Cells.Item(10, 28).Value = 0.354
$ws.Cells.Item(10, 29).Value = 0.706
$ws.Cells.Item(10, 30).Value = 0.208
$ws.Cells.Item(10, 31).Value = 0.456
$ws.Cells.Item(10, 32).Value = 0.971
$ws.Cells.Item(10, 34).Value = 0.169
$ws.Cells.Item(10, 35).Value = 0.794
$ws.Cells.Item(10, 36).Value = 0.163
$ws.Cells.Item(10, 37).Value = 0.404
$ws.Cells.Item(10, 38).Value = 0.9409999999999999
$ws.Cells.Item(10, 39).Value = 0.055
$ws.Cells.Item(10, 40).Value = 0.235
$ws.Cells.Item(10, 41).Value = 0.902
$ws.Cells.Item(11, 2).Value = 0.824
$ws.Cells.Item(11, 3).Value = 0.145
$ws.Cells.Item(11, 4).Value = 0.381
$ws.Cells.Item(11, 5).Value = 0.706
$ws.Cells.Item(11, 6).Value = 0.208
$ws.Cells.Item(11, 7).Value = 0.456
$ws.Cells.Item(11, 8).Value = 0.882
$ws.Cells.Item(11, 9).Value = 0.104
$ws.Cells.Item(11, 10).Value = 0.322
$ws.Cells.Item(11, 11).Value = 0.735
$ws.Cells.Item(11, 12).Value = 0.195
$ws.Cells.Item(11, 13).Value = 0.441
$ws.Cells.Item(11, 14).Value = 0.882
$ws.Cells.Item(11, 15).Value = 0.104
$ws.Cells.Item(11, 16).Value = 0.322
$ws.Cells.Item(11, 17).Value = 0.647
$ws.Cells.Item(11, 18).Value = 0.228
$ws.Cells.Item(11, 19).Value = 0.478
$ws.Cells.Item(11, 20).Value = 0.647
$ws.Cells.Item(11, 21).Value = 0.228
$ws.Cells.Item(11, 22).Value = 0.478
$ws.Cells.Item(11, 23).Value = 0.794
$ws.Cells.Item(11, 24).Value = 0.163
$ws.Cells.Item(11, 25).Value = 0.404
$ws.Cells.Item(11, 26).Value = 0.853
$ws.Cells.Item(11, 27).Value = 0.125
$ws.Cells.Item(11, 28).Value = 0.354
$ws.Cells.Item(11, 29).Value = 0.765
$ws.Cells.Item(11, 30).Value = 0.18
$ws.Cells.Item(11, 31).Value = 0.424
$ws.Cells.Item(11, 32).Value = 0.971
$ws.Cells.Item(11, 34).Value = 0.169
$ws.Cells.Item(11, 35).Value = 0.794
$ws.Cells.Item(11, 36).Value = 0.163
$ws.Cells.Item(11, 37).Value = 0.404
$ws.Cells.Item(11, 38).Value = 0.9409999999999999
$ws.Cells.Item(11, 39).Value = 0.055
$ws.Cells.Item(11, 40).Value = 0.235
$ws.Cells.Item(11, 41).Value = 0.902
$ws.Cells.Item(12, 2).Value = 1.25
$ws.Cells.Item(12, 3).Value = 0.33
$ws.Cells.Item(12, 4).Value = 0.575
$ws.Cells.Item(12, 5).Value = 1.667
$ws.Cells.Item(12, 6).Value = 1.139
$ws.Cells.Item(12, 7).Value = 1.067
$ws.Cells.Item(12, 8).Value = 1.6
$ws.Cells.Item(12, 9).Value = 1.373
$ws.Cells.Item(12, 10).Value = 1.172
$ws.Cells.Item(12, 11).Value = 1.4
$ws.Cells.Item(12, 12).Value = 0.5600000000000001
$ws.Cells.Item(12, 13).Value = 0.748
$ws.Cells.Item(12, 14).Value = 1.367
$ws.Cells.Item(12, 15).Value = 0.5659999999999999
$ws.Cells.Item(12, 16).Value = 0.752
$ws.Cells.Item(12, 26).Value = 1.241
$ws.Cells.Item(12, 27).Value = 0.321
$ws.Cells.Item(12, 28).Value = 0.5669999999999999
$ws.Cells.Item(12, 29).Value = 1.821
$ws.Cells.Item(12, 30).Value = 2.504
$ws.Cells.Item(12, 31).Value = 1.582
$ws.Cells.Item(12, 32).Value = 1.242
$ws.Cells.Item(12, 33).Value = 0.244
$ws.Cells.Item(12, 34).Value = 0.494
$ws.Cells.Item(12, 35).Value = 1.037
$ws.Cells.Item(12, 36).Value = 0.036
$ws.Cells.Item(12, 37).Value = 0.189
$ws.Cells.Item(12, 38).Value = 1.094
$ws.Cells.Item(12, 39).Value = 0.08500000000000001
$ws.Cells.Item(12, 40).Value = 0.291
$ws.Cells.Item(12, 41).Value = 1.124
$ws.Cells.Item(13, 2).Value = 3.441
$ws.Cells.Item(13, 3).Value = 1.423
$ws.Cells.Item(13, 4).Value = 1.193
$ws.Cells.Item(13, 5).Value = 4.571
$ws.Cells.Item(13, 6).Value = 0.459
$ws.Cells.Item(13, 7).Value = 0.678
$ws.Cells.Item(13, 8).Value = 4.594
$ws.Cells.Item(13, 9).Value = 0.679
$ws.Cells.Item(13, 10).Value = 0.824
$ws.Cells.Item(13, 11).Value = 2.265
$ws.Cells.Item(13, 12).Value = 0.606
$ws.Cells.Item(13, 13).Value = 0.779
$ws.Cells.Item(13, 14).Value = 3.235
$ws.Cells.Item(13, 15).Value = 0.768
$ws.Cells.Item(13, 16).Value = 0.876
$ws.Cells.Item(13, 26).Value = 2.515
$ws.Cells.Item(13, 27).Value = 2.916
$ws.Cells.Item(13, 28).Value = 1.708
$ws.Cells.Item(13, 29).Value = 6.353
$ws.Cells.Item(13, 30).Value = 2.228
$ws.Cells.Item(13, 31).Value = 1.493
$ws.Cells.Item(13, 32).Value = 1.588
$ws.Cells.Item(13, 33).Value = 0.595
$ws.Cells.Item(13, 34).Value = 0.771
$ws.Cells.Item(13, 35).Value = 1.206
$ws.Cells.Item(13, 36).Value = 0.163
$ws.Cells.Item(13, 37).Value = 0.404
$ws.Cells.Item(13, 38).Value = 1.5
$ws.Cells.Item(13, 39).Value = 0.721
$ws.Cells.Item(13, 40).Value = 0.849
$ws.Cells.Item(13, 41).Value = 1.431
